# Update "gh-pages" output workbook (江西-漫展信息.xlsx) to match the newly
# scraped data: a handful of "want to go" counts (column F) ticked up, and a
# brand-new event ("南昌·幻梦境国际动漫游戏嘉年华1th", 2024-08-03) was
# appended to both the "展览" sheet and the combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (sheet1): bump column F counters, then append row 30
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

$wsExpo.Range("F5").Value = 5013
$wsExpo.Range("F9").Value = 554
$wsExpo.Range("F11").Value = 1037
$wsExpo.Range("F13").Value = 1400
$wsExpo.Range("F14").Value = 3677
$wsExpo.Range("F17").Value = 123
$wsExpo.Range("F19").Value = 2672
$wsExpo.Range("F22").Value = 87
$wsExpo.Range("F24").Value = 178
$wsExpo.Range("F25").Value = 58
$wsExpo.Range("F28").Value = 270

# New row 30 - copy formatting from the row above (row 29) so column A
# keeps the bold/centered/bordered "index" style, then fill in the values.
$wsExpo.Range("A29:I29").Copy()
$wsExpo.Range("A30:I30").PasteSpecial(-4122)

$wsExpo.Range("A30").Value = 29
# Force column B to stay plain text - otherwise Excel auto-detects the
# "2024-08-03" string as a date and stores a date serial instead.
$wsExpo.Range("B30").NumberFormat = "@"
$wsExpo.Range("B30").Value = "2024-08-03"
$wsExpo.Range("C30").Value = "南昌·幻梦境国际动漫游戏嘉年华1th"
$wsExpo.Range("D30").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
$wsExpo.Range("E30").Value = "2024.08.03 09:00-08.04 17:30"
$wsExpo.Range("F30").Value = 1
$wsExpo.Range("G30").Value = 64
$wsExpo.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=83980"
$wsExpo.Range("I30").Value = "//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg"

# ---------------------------------------------------------------------
# Sheet "演出" (sheet2): bump the single counter
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 43

# ---------------------------------------------------------------------
# Sheet "全部类型" (sheet4): same counter bumps (shifted by the extra
# "演出" row at position 5), then append row 31 with the same new event.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F5").Value = 43
$wsAll.Range("F6").Value = 5013
$wsAll.Range("F10").Value = 554
$wsAll.Range("F12").Value = 1037
$wsAll.Range("F14").Value = 1400
$wsAll.Range("F15").Value = 3677
$wsAll.Range("F18").Value = 123
$wsAll.Range("F20").Value = 2672
$wsAll.Range("F23").Value = 87
$wsAll.Range("F25").Value = 178
$wsAll.Range("F26").Value = 58
$wsAll.Range("F29").Value = 270

$wsAll.Range("A30:I30").Copy()
$wsAll.Range("A31:I31").PasteSpecial(-4122)

$wsAll.Range("A31").Value = 30
$wsAll.Range("B31").NumberFormat = "@"
$wsAll.Range("B31").Value = "2024-08-03"
$wsAll.Range("C31").Value = "南昌·幻梦境国际动漫游戏嘉年华1th"
$wsAll.Range("D31").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
$wsAll.Range("E31").Value = "2024.08.03 09:00-08.04 17:30"
$wsAll.Range("F31").Value = 1
$wsAll.Range("G31").Value = 64
$wsAll.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=83980"
$wsAll.Range("I31").Value = "//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg"
